$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","G","J","K","N","O")

$newValues = @{
    "B" = @(0.6261066896930174, 0.5870295664328182, 0.5633033008686539, 0.553702232698015, 0.5521120753094806, 0.5631735433862275, 0.6125776428759764, 0.711568362696795, 0.7855760813206984, 0.8195209748376442, 0.8324147990271058, 0.8296361271918329, 0.8205809655947576, 0.8150395665949191, 0.7833632703606384, 0.7640019075961106, 0.7528919933730549, 0.7491348846322126, 0.766060245995078, 0.8232396137438798, 0.860840440801212, 0.8407511868276174, 0.7651296051659244, 0.684563620834183)
    "C" = @(0.05682355408124806, 0.05105316440096885, 0.04751945293608628, 0.04608182339397615, 0.04584325144101342, 0.04750005484477526, 0.05483200926600773, 0.06928296945915235, 0.07994434051781241, 0.08480409905259023, 0.08664575554249154, 0.08624906158073031, 0.08495558601059372, 0.08416347215316478, 0.07962693820168454, 0.07684641186315844, 0.07524805618911046, 0.07470704143251794, 0.07714230776056752, 0.0853354742751975, 0.09069815529014136, 0.0878352766477235, 0.07700853258110385, 0.06536581291118182)
    "D" = @(0.6622854274184249, 0.6521792814765206, 0.6463117072544833, 0.6440056015984226, 0.6436278092765804, 0.6462802621294657, 0.6587307917062333, 0.6858240307079768, 0.7073633476507553, 0.7175174539377167, 0.7214136860673648, 0.7205722914575006, 0.717836975706291, 0.7161681676285241, 0.7067069061459961, 0.7009938011076997, 0.6977412641982426, 0.6966457650415805, 0.7015985052820497, 0.7186390181356046, 0.7300737540943487, 0.7239435954250837, 0.7013250188350355, 0.6782077670391686)
    "E" = @(0.2700698995109008, 0.2650986388462329, 0.2621903322277888, 0.2610414130742171, 0.2608528250792403, 0.2621746907528362, 0.268325912764773, 0.28153199162454, 0.2919336926731901, 0.2968180552374236, 0.2986895900560214, 0.2982855462025569, 0.2969715880118713, 0.2961696067455506, 0.2916175576691344, 0.2888640960531106, 0.2872947394345786, 0.2867658491149285, 0.2891557204195365, 0.2973569341952711, 0.3028447387794913, 0.2999040995339115, 0.2890238345412826, 0.2778367607061085)
    "G" = @(0.002503748502152242, 0.002506794515150612, 0.002508764963276613, 0.00250959320581325, 0.002509732263470489, 0.002508776030650498, 0.002504778020931833, 0.002497729320533128, 0.002493028147587828, 0.002490992097054869, 0.002490235762553074, 0.002490398001079783, 0.002490929579114102, 0.002491257094833554, 0.00249316326566329, 0.002494358851005158, 0.002495056175667357, 0.002495293938195142, 0.002494230580206423, 0.002490773044220296, 0.002488598842076576, 0.002489751454947831, 0.002494288540448901, 0.002499551970536264)
    "J" = @(0.1405152317203644, 0.1372110889403189, 0.1352617626820418, 0.1344873593653659, 0.1343599757007325, 0.1352512379919659, 0.1393594696580465, 0.1480470305573363, 0.1548170347979436, 0.1579815402118214, 0.1591920770732429, 0.1589308231573199, 0.1580808870754709, 0.1575618668488374, 0.1546119344193357, 0.1528239745688467, 0.151803567981446, 0.1514594461079071, 0.1530134798840521, 0.1583302024919391, 0.1618761482147022, 0.159977094536444, 0.1529277811237932, 0.1456289814951361)
    "K" = @(0.6056669903842078, 0.5634197884708101, 0.5377146719908126, 0.5272990582394925, 0.5255731541371063, 0.5375739622495814, 0.5910516274083477, 0.6977736828757202, 0.7773058050520092, 0.8137302660756518, 0.8275582249521847, 0.824578585107929, 0.8148672053567338, 0.8089232254608305, 0.7749302737711332, 0.7541391313415033, 0.7422036916946695, 0.7381665326825555, 0.7563499986871989, 0.8177187316360346, 0.8580294326811497, 0.8364964423145409, 0.7553504105050308, 0.6687048311693502)
    "N" = @(2.041400171021401, 2.06170018011521, 2.074808703301706, 2.080312485781562, 2.08123616494332, 2.074882273560178, 2.048265877539631, 2.001184176447818, 1.969710842547229, 1.95606990719887, 1.951001772181915, 1.95208895229645, 1.955650996119765, 1.95784553762671, 1.970615931637859, 1.97862349174158, 1.983292891052447, 1.984884799951995, 1.977764483604346, 1.954602093364283, 1.940031860599785, 1.947756280423388, 1.978152636385683, 2.013373598746249)
    "O" = @(5.156456206425617, 5.167066794255078, 5.176082454028034, 5.180385133652692, 5.181137561009535, 5.176137935931223, 5.159595641771801, 5.147007922203301, 5.149882498298467, 5.153827624869592, 5.155701109199214, 5.155280735602133, 5.153974149139231, 5.15322326446136, 5.149677771902105, 5.148178447476425, 5.147564309157104, 5.147398998204864, 5.14831235982291, 5.154347622056974, 5.160504462290049, 5.157015905826967, 5.148251046032669, 5.148285862603359)
}

$startRow = 2
foreach ($col in $cols) {
    $vals = $newValues[$col]
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Range($col + ($startRow + $i)).Value = $vals[$i]
    }
}
